$wb = $excel.ActiveWorkbook

# --- Measures sheet: insert new "item_num" column before the old column O ("comment") ---
$ws = $wb.Worksheets.Item("Measures")
$ws.Columns("O").Insert()
$ws.Range("O1").Value = "item_num"
$ws.Range("O2").Value = 1
$ws.Range("O3").Value = 1
$ws.Range("O4").Value = 1
$ws.Range("O5").Select()

# --- ID sheet: same new column, no data values in the new column ---
$ws = $wb.Worksheets.Item("ID")
$ws.Columns("O").Insert()
$ws.Range("O1").Value = "item_num"
$ws.Range("O2").Select()

# The hidden AutoFilter-derived named range needs to grow by one column too
foreach ($n in $wb.Names) {
    if ($n.Name -eq "ID!_FilterDatabase") {
        $n.RefersTo = "=ID!`$A`$1:`$P`$1"
    }
}

# --- Dems sheet ---
$ws = $wb.Worksheets.Item("Dems")
$ws.Columns("O").Insert()
$ws.Range("O1").Value = "item_num"
$ws.Range("O2").Select()

# --- Dates sheet ---
$ws = $wb.Worksheets.Item("Dates")
$ws.Columns("O").Insert()
$ws.Range("O1").Value = "item_num"
$ws.Range("O2").Select()

# --- NewVars sheet ---
$ws = $wb.Worksheets.Item("NewVars")
$ws.Columns("O").Insert()
$ws.Range("O1").Value = "item_num"
$ws.Range("O2").Select()

# Restore the originally active sheet/tab (PanelInfo) which was not otherwise touched
$wb.Worksheets.Item("PanelInfo").Activate()
